$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate performance numbers
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1000.15
$summary.Range("B4").Value = 0.16
$summary.Range("B5").Value = 0.32
$summary.Range("B6").Value = 10
$summary.Range("B7").Value = 5
$summary.Range("B9").Value = 50

# ---------------------------------------------------------------------------
# Strategy Status sheet: leadlag strategy row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.13
$status.Range("D5").Value = 9
$status.Range("E5").Value = 0.14
$status.Range("F5").Value = 0.13
$status.Range("G5").Value = 44.44

# ---------------------------------------------------------------------------
# All Trades sheet: trade #10 (row 11) moves from OPEN to CLOSED with real
# Polymarket fill/slippage data
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G11").Value = 69157.66054500001
$allTrades.Range("H11").Value = "CLOSED"
$allTrades.Range("I11").Value = 1.202
$allTrades.Range("J11").Value = 0.06
$allTrades.Range("K11").Value = 100.13
$allTrades.Range("N11").Value = "time_exit_5min"
$allTrades.Range("O11").Value = 5

# ---------------------------------------------------------------------------
# leadlag sheet: same trade, mirrored on the strategy-specific tab (row 10)
# ---------------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")
$leadlag.Range("G10").Value = 69157.66054500001
$leadlag.Range("H10").Value = "CLOSED"
$leadlag.Range("I10").Value = 1.202
$leadlag.Range("J10").Value = 0.06
$leadlag.Range("K10").Value = 100.13
$leadlag.Range("N10").Value = "time_exit_5min"
$leadlag.Range("O10").Value = 5
